$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 5 - this shifts the existing rows 5..21 down to 6..22,
# carrying all their data (and the date-format style on column D) along with them.
$ws.Rows.Item(5).Insert()

# Populate the newly inserted row 5 with the new weekly record.
$ws.Range("A5").Value = 5
$ws.Range("B5").Value = "Macroferia Regional de Talca"
$ws.Range("C5").Value = "Maule"
$ws.Range("D5").Value = 44469
$ws.Range("E5").Value = 7
$ws.Range("F5").Value = 300000000
$ws.Range("G5").Value = "Espárragos"
$ws.Range("H5").Value = "Verde"
$ws.Range("I5").Value = "Primera"
$ws.Range("J5").Value = 3000
$ws.Range("K5").Value = 1200
$ws.Range("L5").Value = 1200
$ws.Range("M5").Value = 1200
$ws.Range("N5").Value = "$/kilo"
$ws.Range("O5").Value = "Provincia de Linares"
$ws.Range("P5").Value = 1200
$ws.Range("Q5").Value = 1
$ws.Range("R5").Value = "Hortaliza"
